# Correção nos dados e início da análise PNAD 2009
#
# The row "grandes regiões e unidades da federação" (row 6) was a blank
# header/divider row with no data values. It is removed, so every row
# below it (7..37) shifts up by one, and the now-unused shared string
# is dropped automatically when the workbook is saved. This also removes
# the last row (37) since the sheet simply has one fewer row overall.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6").Delete()
